$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trafos")

# Insert a new "v_base_kV" column between "V_lv_kV" (D) and "V_SCH_pu" (old E)
$ws.Columns("E:E").Insert()
$ws.Range("E1").Value = "v_base_kV"
$ws.Range("E2").Value = 132

# Make the trafos sheet the active sheet/tab (was "loads" before)
$ws.Select()

# Select the full header + data rows (rows 1:2), matching the new working selection
$ws.Rows("1:2").Select()
